# Proj_Data_CSV.xlsx -- "Fixed Book table and advanced queries"
#
# While reviewing/updating the advanced queries (Transaction_id helper
# column) the author scrolled further down the sheet and selected the
# whole of column K. Re-touching the workbook this way also makes Excel
# recalculate every volatile formula on the sheet -- the O/R/S helper
# columns are driven by RANDBETWEEN(), used to randomly assign
# shop/shelf/warehouse placement for the generated INSERT statements --
# so those cached results move around too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Scroll down to roughly where the author ended up (view anchored a few
# rows further down than before) and select the entire Transaction_id
# column (K), matching the saved selection sqref "K1:K1048576".
$ws.Range("A83").Select() | Out-Null
$ws.Columns("K").Select() | Out-Null

# RANDBETWEEN() in columns O, R and S is volatile; force a full workbook
# recalculation so every one of those cells gets a freshly generated
# cached value, same as happens whenever the workbook is touched/saved.
$excel.CalculateFull() | Out-Null
